$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 9
$ws.Range("F6").Value = 1
$ws.Range("F8").Value = 6
$ws.Range("F10").Value = -1
$ws.Range("F15").Value = 5
$ws.Range("F18").Value = 3
